$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 16: drawing
$ws.Range("D15").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = "drawing:null (demanar per celeste)"

# New row 17: phrase
$ws.Range("D15").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").Value = "phrase:string(conté la frase)"

# Rename the "tournId:Number" entry (row 15, col E) to the new, longer description
$ws.Range("E15").Value = "turnId:Number(demanat per Oriol)"

# Widen column E to fit the new text
$ws.Columns.Item(5).ColumnWidth = 46.25

# Update view: zoom + selection (mirrors the author re-inspecting the sheet)
$ws.Application.ActiveWindow.Zoom = 154
$ws.Range("E26").Select()
